$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.937.25"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.638.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "214.60"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "  -0.23%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("E8").Value = "  -0.38%  "
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("E10").Value = "  -0.68%  "
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.25"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.636.00"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("E14").Value = "  -0.48%  "
$ws.Range("E15").Value = "  +0.88%  "
$ws.Range("E16").Value = "  -0.45%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "25.980.51"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "193.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.22%  "
$ws.Range("E20").Value = "  -0.98%  "
$ws.Range("E21").Value = "  -0.73%  "
$ws.Range("E22").Value = "  -1.71%  "
$ws.Range("E23").Value = "  +3.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "143.63"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.41%  "
$ws.Range("E25").Value = "  -0.20%  "
$ws.Range("E26").Value = "  -0.66%  "
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  -0.13%  "
$ws.Range("E29").Value = "  -0.08%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0493"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.64%  "
$ws.Range("E31").Value = "  -0.85%  "
$ws.Range("E32").Value = "  +0.34%  "
$ws.Range("E34").Value = "  +0.67%  "
$ws.Range("E35").Value = "  -0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.129.08"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.97%  "
$ws.Range("E37").Value = "  -1.34%  "
$ws.Range("E38").Value = "  -0.16%  "
$ws.Range("E39").Value = "  -0.71%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.63"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.791"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("E43").Value = "  -0.31%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "56.31"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.48%  "
$ws.Range("E45").Value = "  +2.39%  "
$ws.Range("E46").Value = "  -1.57%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "7.72"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.62%  "
$ws.Range("E48").Value = "  -0.48%  "
$ws.Range("E50").Value = "  -2.08%  "
$ws.Range("E51").Value = "  -0.65%  "
